{"js": "// The interview transcript alternates speaker labels \"Interviewer\" / \"Interviewee\"\n// at the start of each paragraph. Two adjacent paragraphs were mislabeled\n// (two \"Interviewer\" rows in a row, followed by two \"Interviewee\" rows in a\n// row) breaking the strict alternation. Correct the mislabeling by swapping\n// the speaker label in those two paragraphs:\n//   \"Interviewer: The fact that it makes people without searching for\n//    firewood is very crucial...\"  ->  should be the Interviewee's answer\n//   \"Interviewee: Why do you think it failed?\"  ->  should be the\n//    Interviewer's question\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetMislabeledInterviewer = null; // paragraph that says \"Interviewer\" but should say \"Interviewee\"\nlet targetMislabeledInterviewee = null; // paragraph that says \"Interviewee\" but should say \"Interviewer\"\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length - 1; i++) {\n  const text = items[i].text;\n  const nextText = items[i + 1].text;\n\n  // Two consecutive paragraphs labeled \"Interviewer\" -> the second one is\n  // mislabeled and should be \"Interviewee\".\n  if (\n    text.indexOf(\"Interviewer:\") === 0 &&\n    nextText.indexOf(\"Interviewer:\") === 0\n  ) {\n    targetMislabeledInterviewer = items[i + 1];\n  }\n\n  // Two consecutive paragraphs labeled \"Interviewee\" -> the first one is\n  // mislabeled and should be \"Interviewer\".\n  if (\n    text.indexOf(\"Interviewee:\") === 0 &&\n    nextText.indexOf(\"Interviewee:\") === 0\n  ) {\n    targetMislabeledInterviewee = items[i];\n  }\n}\n\nif (!targetMislabeledInterviewer || !targetMislabeledInterviewee) {\n  throw new Error(\"Could not locate the mislabeled Interviewer/Interviewee paragraphs.\");\n}\n\n// Replace just the speaker-label run's text, leaving the rest of each\n// paragraph (and all character formatting) untouched.\nconst interviewerLabelRange = targetMislabeledInterviewer.search(\"Interviewer\", {\n  matchCase: true,\n  matchWholeWord: true,\n}).getFirst();\nconst intervieweeLabelRange = targetMislabeledInterviewee.search(\"Interviewee\", {\n  matchCase: true,\n  matchWholeWord: true,\n}).getFirst();\n\ninterviewerLabelRange.insertText(\"Interviewee\", Word.InsertLocation.replace);\nintervieweeLabelRange.insertText(\"Interviewer\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The interview transcript alternates speaker labels \"Interviewer\" / \"Interviewee\"\n# at the start of each paragraph. Two adjacent paragraphs were mislabeled\n# (two \"Interviewer\" rows in a row, followed by two \"Interviewee\" rows in a\n# row) breaking the strict alternation:\n#   \"Interviewer: The fact that it makes people without searching for\n#    firewood is very crucial...\"  -> should be the Interviewee's answer\n#   \"Interviewee: Why do you think it failed?\"  -> should be the\n#    Interviewer's question\n# Correct the mislabeling by swapping just the bold speaker-label word in\n# those two paragraphs, leaving everything else untouched.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n$mislabeledInterviewerIndex = 0\n$mislabeledIntervieweeIndex = 0\n\nfor ($i = 1; $i -lt $count; $i++) {\n  $textCur = $d.Paragraphs($i).Range.Text\n  $textNext = $d.Paragraphs($i + 1).Range.Text\n\n  # Two consecutive paragraphs labeled \"Interviewer:\" -> the second one is\n  # mislabeled and should be \"Interviewee:\".\n  if ($textCur.StartsWith(\"Interviewer:\") -and $textNext.StartsWith(\"Interviewer:\")) {\n    $mislabeledInterviewerIndex = $i + 1\n  }\n\n  # Two consecutive paragraphs labeled \"Interviewee:\" -> the first one is\n  # mislabeled and should be \"Interviewer:\".\n  if ($textCur.StartsWith(\"Interviewee:\") -and $textNext.StartsWith(\"Interviewee:\")) {\n    $mislabeledIntervieweeIndex = $i\n  }\n}\n\nif ($mislabeledInterviewerIndex -eq 0 -or $mislabeledIntervieweeIndex -eq 0) {\n  throw \"Could not locate the mislabeled Interviewer/Interviewee paragraphs.\"\n}\n\n# Replace just the speaker-label word, leaving the rest of each paragraph\n# (and all character formatting) untouched.\n$labelRangeToFixInterviewer = $d.Paragraphs($mislabeledInterviewerIndex).Range.Words(1)\n$labelRangeToFixInterviewee = $d.Paragraphs($mislabeledIntervieweeIndex).Range.Words(1)\n\n$labelRangeToFixInterviewer.Text = \"Interviewee\"\n$labelRangeToFixInterviewee.Text = \"Interviewer\"\n"}
